$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The three observation rows (2, 3, 4) get cyclically rotated:
#   old row 4 -> row 2
#   old row 2 -> row 3
#   old row 3 -> row 4
# Row 3's species/location data already matches row 2's (same species block),
# so only its Id/Ost/Nord need updating. Row 2 and row 4 swap species records
# (a woodpecker record moves out, a lichen record moves in, and vice versa),
# so every species-specific column is rewritten explicitly below.
# ---------------------------------------------------------------------------

# --- Row 2: becomes the "Garnlav" (lichen) record that used to be row 4 ---
$ws.Range("A2").Value2 = 131223489
$ws.Range("B2").Value2 = 79245
$ws.Range("E2").Value2 = 6425
$ws.Range("F2").Value2 = "Garnlav"
$ws.Range("G2").Value2 = "Alectoria sarmentosa"
$ws.Range("H2").Value2 = "(Ach.) Ach."
$ws.Range("J2").Value2 = "bålar"
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("Q2").Value2 = 773011
$ws.Range("R2").Value2 = 7122664
$ws.Range("AC2").ClearContents()
$ws.Range("AF2").Value2 = ""

# --- Row 3: stays the "Tretåig hackspett" record, only Id/coords shift ---
$ws.Range("A3").Value2 = 131223149
$ws.Range("Q3").Value2 = 772974
$ws.Range("R3").Value2 = 7122563

# --- Row 4: becomes the "Tretåig hackspett" record that used to be row 3 ---
$ws.Range("A4").Value2 = 131223060
$ws.Range("B4").Value2 = 57884
$ws.Range("E4").Value2 = 100109
$ws.Range("F4").Value2 = "Tretåig hackspett"
$ws.Range("G4").Value2 = "Picoides tridactylus"
$ws.Range("H4").Value2 = "(Linnaeus, 1758)"
$ws.Range("J4").ClearContents()
$ws.Range("L4").Value2 = ""
$ws.Range("M4").Value2 = "färska spår"
$ws.Range("Q4").Value2 = 772981
$ws.Range("R4").Value2 = 7122639
$ws.Range("AC4").Value2 = "färska ringhack på gran"
$ws.Range("AF4").ClearContents()

# --- Row 16: Taxonsorteringsordning value correction ---
$ws.Range("B16").Value2 = 91813
